$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A7: long whitespace-padded text ending in "aaa" (keeps style s="3")
$ws.Range("A7").Value = "                                                                                                                                                                                                                                                                                                    aaa"

# B7: email text, loses its explicit style (falls back to default style 0)
$ws.Range("B7").Value = "ahshshssh@gmail.com"
$ws.Range("B7").ClearFormats()

# C7 / D7: value "111111" (same shared string used by row 2/3), with the
# quote-prefixed style (s="5") instead of the plain style (s="3")
$ws.Range("C7").Value = "111111"
$ws.Range("D7").Value = "111111"
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)

# M7: boolean TRUE
$ws.Range("M7").Value = $true

# Row 7 height grows
$ws.Rows.Item(7).RowHeight = 23.25

# Selection moves to A7
$ws.Range("A7").Select()
